$d = $word.ActiveDocument
$tbl = $d.Tables.Item(3)

# --- Right-align the price column paragraphs (ContainerPrice, DiscountValue,
#     AdditionalCostPrice, PriceWithoutTax, TaxValue, TotalPrice) ---
$priceRows = 1,2,3,4,5,6
foreach ($r in $priceRows) {
    $p = $tbl.Cell($r, 3).Range.Paragraphs.Item(1)
    $p.Alignment = 2   # wdAlignParagraphRight
}

Write-Host "Done alignments"
